$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") holds the "last changed" date for every
# record, stored as the serial date number 45171 (2023-09-02). The update
# bumps this to 45172 (2023-09-03) for every data row (rows 2 through 527).
$ws.Range("C2:C527").Value = 45172
